$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "https://www.br.de/index.html"
$ws.Range("A5").Value = "testing 7 / 7"
